# Update "想去人数" (want-to-go count) figures on the 展览 (Exhibition) and
# 全部类型 (All Types) sheets to match the latest scrape output.

$wb = $excel.ActiveWorkbook

# Sheet "展览" - row => new value for column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 1333
$wsExpo.Range("F3").Value = 1220
$wsExpo.Range("F4").Value = 14625
$wsExpo.Range("F5").Value = 17780
$wsExpo.Range("F7").Value = 73
$wsExpo.Range("F9").Value = 219
$wsExpo.Range("F24").Value = 7338
$wsExpo.Range("F26").Value = 5
$wsExpo.Range("F28").Value = 1178
$wsExpo.Range("F30").Value = 5873
$wsExpo.Range("F31").Value = 70
$wsExpo.Range("F35").Value = 229
$wsExpo.Range("F36").Value = 5110

# Sheet "全部类型" - row => new value for column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 1333
$wsAll.Range("F3").Value = 1220
$wsAll.Range("F4").Value = 14625
$wsAll.Range("F5").Value = 17780
$wsAll.Range("F7").Value = 73
$wsAll.Range("F9").Value = 219
$wsAll.Range("F25").Value = 7338
$wsAll.Range("F27").Value = 5
$wsAll.Range("F29").Value = 1178
$wsAll.Range("F32").Value = 5873
$wsAll.Range("F33").Value = 70
$wsAll.Range("F37").Value = 229
$wsAll.Range("F38").Value = 5110

$wb.Save()
